$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new log entry row (row 23): hours and description for the new task.
$ws.Range("B23").Value() = 1
$ws.Range("C23").Value() = "Reduce GPU computation when loading GIF."

# Extend the total-hours sum formula to include the new row.
$ws.Range("B27").Formula() = "=SUM(B2:B23)"

# Move the active selection to C24, as in the saved workbook.
[void]$ws.Range("C24").Select()
